# Apply the Sep 19 2024 11:30:04 UTC cryptos-list refresh (GitHub Actions bot).
# Updates Price (D) / Volume 1h % (E) figures, and fixes the Fetch.AI /
# Binance-PegBSC-USD row ordering (rows 27-28 swapped their Coin/Link pair).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "62.499.00" }
    @{ Cell = "E2"; Value = "  +4.20%  " }
    @{ Cell = "D3"; Value = "2.427.52" }
    @{ Cell = "E3"; Value = "  +5.18%  " }
    @{ Cell = "E4"; Value = "  +0.02%  " }
    @{ Cell = "D5"; Value = "557.09" }
    @{ Cell = "E5"; Value = "  +2.82%  " }
    @{ Cell = "D6"; Value = "138.87" }
    @{ Cell = "E6"; Value = "  +7.22%  " }
    @{ Cell = "E7"; Value = "  +0.02%  " }
    @{ Cell = "E8"; Value = "  +2.13%  " }
    @{ Cell = "D9"; Value = "2.425.14" }
    @{ Cell = "E10"; Value = "  +3.73%  " }
    @{ Cell = "E11"; Value = "  +3.91%  " }
    @{ Cell = "E12"; Value = "  +0.39%  " }
    @{ Cell = "D13"; Value = "0.349" }
    @{ Cell = "E13"; Value = "  +4.97%  " }
    @{ Cell = "D14"; Value = "26.19" }
    @{ Cell = "E14"; Value = "  +12.13%  " }
    @{ Cell = "D15"; Value = "2.860.48" }
    @{ Cell = "E15"; Value = "  +5.20%  " }
    @{ Cell = "D16"; Value = "62.370.60" }
    @{ Cell = "E16"; Value = "  +4.07%  " }
    @{ Cell = "E17"; Value = "  +7.31%  " }
    @{ Cell = "D18"; Value = "2.431.39" }
    @{ Cell = "E18"; Value = "  +3.97%  " }
    @{ Cell = "E19"; Value = "  +6.52%  " }
    @{ Cell = "D20"; Value = "345.70" }
    @{ Cell = "E20"; Value = "  +10.80%  " }
    @{ Cell = "E21"; Value = "  +3.23%  " }
    @{ Cell = "D22"; Value = "6.84" }
    @{ Cell = "E22"; Value = "  +4.18%  " }
    @{ Cell = "E23"; Value = "  -0.03%  " }
    @{ Cell = "D24"; Value = "5.55" }
    @{ Cell = "E24"; Value = "  -2.59%  " }
    @{ Cell = "E26"; Value = "  +1.30%  " }
    @{ Cell = "B27"; Value = "Binance-PegBSC-USD" }
    @{ Cell = "C27"; Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd" }
    @{ Cell = "D27"; Value = "1.00" }
    @{ Cell = "E27"; Value = "  +0.30%  " }
    @{ Cell = "B28"; Value = "Fetch.AI" }
    @{ Cell = "C28"; Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet" }
    @{ Cell = "D28"; Value = "1.54" }
    @{ Cell = "E28"; Value = "  +14.37%  " }
    @{ Cell = "E29"; Value = "  +5.20%  " }
    @{ Cell = "D30"; Value = "1.35" }
    @{ Cell = "E30"; Value = "  +15.23%  " }
    @{ Cell = "E31"; Value = "  +5.34%  " }
    @{ Cell = "E32"; Value = "  +8.32%  " }
    @{ Cell = "D33"; Value = "6.47" }
    @{ Cell = "E33"; Value = "  +11.13%  " }
    @{ Cell = "D34"; Value = "172.19" }
    @{ Cell = "E34"; Value = "  +0.44%  " }
    @{ Cell = "D35"; Value = "1.44" }
    @{ Cell = "E35"; Value = "  +5.79%  " }
    @{ Cell = "E36"; Value = "  +4.58%  " }
    @{ Cell = "D37"; Value = "379.94" }
    @{ Cell = "E37"; Value = "  +19.99%  " }
    @{ Cell = "D38"; Value = "18.57" }
    @{ Cell = "D39"; Value = "4.46" }
    @{ Cell = "E39"; Value = "  +11.60%  " }
    @{ Cell = "E41"; Value = "  -0.04%  " }
    @{ Cell = "E42"; Value = "  +11.80%  " }
    @{ Cell = "D43"; Value = "39.34" }
    @{ Cell = "E43"; Value = "  +3.28%  " }
    @{ Cell = "D44"; Value = "145.02" }
    @{ Cell = "E44"; Value = "  +6.59%  " }
    @{ Cell = "D45"; Value = "3.67" }
    @{ Cell = "E45"; Value = "  +7.31%  " }
    @{ Cell = "D46"; Value = "20.80" }
    @{ Cell = "E46"; Value = "  +10.64%  " }
    @{ Cell = "D47"; Value = "0.591" }
    @{ Cell = "E47"; Value = "  +4.12%  " }
    @{ Cell = "D48"; Value = "0.0954" }
    @{ Cell = "E48"; Value = "  +1.87%  " }
    @{ Cell = "D49"; Value = "0.0520" }
    @{ Cell = "E49"; Value = "  +6.13%  " }
    @{ Cell = "E50"; Value = "  +4.85%  " }
    @{ Cell = "D51"; Value = "17.86" }
    @{ Cell = "E51"; Value = "  +6.26%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Cell.StartsWith("D")) {
        # Keep price column as text so values like "1.00" / "345.70" / "62.499.00"
        # do not get reinterpreted as numbers and lose their formatting.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
